$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Patient row 3 (Shehab) - values entered through the form / manual input
$ws.Range("B3").Value = "Shehab"
$ws.Range("M3").Value = "O+"
$ws.Range("O3").Value = 25
$ws.Range("P3").Value = 42502.77497329861

Write-Host "Updated patient row 3"
